$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.323.04'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.27%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.934.41'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.41%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7505'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.26%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.72'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.34%  '

$ws.Range('E7').Value = '  +0.10%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '27.99'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.39%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3179'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.81%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07229'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.33%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7797'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.56%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08041'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.98%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.932.63'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.07%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.399'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.32%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.98'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.83%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.50'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.13%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.333.68'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.23%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.111'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +6.07%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '251.84'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.32%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000008022'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.25%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.226.94'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.91%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9997'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.09%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.17%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.695'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.35%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.558'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.18%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.89'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.60%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.10'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.02%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1305'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.55%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.198'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.63%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.381'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.91%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.543'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.16%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.430'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.97%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.150'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.71%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05290'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.02%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.332'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.09%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7561'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.69%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.790'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.86%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01957'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.08%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.801'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.24%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '78.91'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.94%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.502'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.09%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4513'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.54%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.982'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.50%  '

$ws.Range('E44').Value = '  +0.12%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8402'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.50%  '

$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.06'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.31%  '

$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.703'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.81%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '101.70'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.10%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '37.62'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.92%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1231'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +8.27%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '964.28'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.24%  '
